$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 32
$ws.Range("H32").Value = 3600.2856
$ws.Range("I32").Value = 1899
$ws.Range("J32").Value = 3883.8333
$ws.Range("K32").Value = 1899
$ws.Range("L32").Value = 3883.8333
$ws.Range("M32").Value = -1573
$ws.Range("N32").Value = -4535.8333
# row 38
$ws.Range("H38").Value = 1592.25
$ws.Range("I38").Value = 1009.5
$ws.Range("J38").Value = 2175
$ws.Range("K38").Value = 3028.5
$ws.Range("L38").Value = 6525
$ws.Range("M38").Value = -2656.5
$ws.Range("N38").Value = -7269
# row 41
$ws.Range("H41").Value = 327.69232
$ws.Range("I41").Value = 194.28572
$ws.Range("J41").Value = 483.33334
$ws.Range("K41").Value = 194.28572
$ws.Range("L41").Value = 483.33334
$ws.Range("M41").Value = 245.71428
$ws.Range("N41").Value = -1363.33334
# row 42
$ws.Range("H42").Value = 108
$ws.Range("I42").Value = 108
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 324
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -94
$ws.Range("N42").ClearContents()
# row 51
$ws.Range("H51").Value = 2635.7144
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 2658.3333
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 2658.3333
$ws.Range("M51").Value = -2016
$ws.Range("N51").Value = -3626.3333
# row 55
$ws.Range("H55").Value = 366.41666
$ws.Range("I55").Value = 299.44446
$ws.Range("J55").Value = 567.3333
$ws.Range("K55").Value = 299.44446
$ws.Range("L55").Value = 567.3333
$ws.Range("M55").Value = -85.44445999999999
$ws.Range("N55").Value = -995.3333
# row 112
$ws.Range("H112").Value = 1311.7838
$ws.Range("I112").Value = 699.6667
$ws.Range("J112").Value = 1365.7941
$ws.Range("K112").Value = 2099.0001
$ws.Range("L112").Value = 4097.3823
$ws.Range("M112").Value = -991.0001000000002
$ws.Range("N112").Value = -6313.3823
# row 127
$ws.Range("H127").Value = 2751.6155
$ws.Range("I127").Value = 433
$ws.Range("J127").Value = 3447.2
$ws.Range("K127").Value = 1299
$ws.Range("L127").Value = 10341.6
$ws.Range("M127").Value = 3661
$ws.Range("N127").Value = -20261.6
# row 129
$ws.Range("H129").Value = 1106.5253
$ws.Range("I129").Value = 598.5
$ws.Range("J129").Value = 1117
$ws.Range("K129").Value = 1795.5
$ws.Range("L129").Value = 3351
$ws.Range("M129").Value = 3204.5
$ws.Range("N129").Value = -13351
# row 138
$ws.Range("H138").Value = 3260.1492
$ws.Range("I138").Value = 994.875
$ws.Range("J138").Value = 4524.4883
$ws.Range("K138").Value = 2984.625
$ws.Range("L138").Value = 13573.4649
$ws.Range("M138").Value = 2155.375
$ws.Range("N138").Value = -23853.4649

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 17245714
$ws.Range("I32").Value = 19610688
$ws.Range("J32").Value = 15175.429
$ws.Range("K32").Value = 19610688
$ws.Range("L32").Value = 15175.429
$ws.Range("M32").Value = -19610401
$ws.Range("N32").Value = -15749.429
# row 61
$ws.Range("H61").Value = 1482.2572
$ws.Range("I61").Value = 1015.5517
$ws.Range("J61").Value = 3738
$ws.Range("K61").Value = 1015.5517
$ws.Range("L61").Value = 3738
$ws.Range("M61").Value = -803.5517
$ws.Range("N61").Value = -4162
# row 74
$ws.Range("H74").Value = 1277.9048
$ws.Range("I74").Value = 1190.7567
$ws.Range("J74").Value = 1922.8
$ws.Range("K74").Value = 1190.7567
$ws.Range("L74").Value = 1922.8
$ws.Range("M74").Value = -316.7566999999999
$ws.Range("N74").Value = -3670.8
# row 77
$ws.Range("H77").Value = 1277.9048
$ws.Range("I77").Value = 1190.7567
$ws.Range("J77").Value = 1922.8
$ws.Range("K77").Value = 5953.7835
$ws.Range("L77").Value = 9614
$ws.Range("M77").Value = -1585.7835
$ws.Range("N77").Value = -18350
# row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# row 136
$ws.Range("H136").Value = 1482.2572
$ws.Range("I136").Value = 1015.5517
$ws.Range("J136").Value = 3738
$ws.Range("K136").Value = 3046.6551
$ws.Range("L136").Value = 11214
$ws.Range("M136").Value = -496.6550999999999
$ws.Range("N136").Value = -16314

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 107
$ws.Range("H107").Value = 3722.5
$ws.Range("I107").Value = 3968.5715
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 3968.5715
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -2048.5715
$ws.Range("N107").Value = -5840

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2271.125
$ws.Range("I31").Value = 1708.6786
$ws.Range("J31").Value = 3583.5
$ws.Range("K31").Value = 1708.6786
$ws.Range("L31").Value = 3583.5
$ws.Range("M31").Value = -1413.6786
$ws.Range("N31").Value = -4173.5
# row 34
$ws.Range("H34").Value = 2271.125
$ws.Range("I34").Value = 1708.6786
$ws.Range("J34").Value = 3583.5
$ws.Range("K34").Value = 1708.6786
$ws.Range("L34").Value = 3583.5
$ws.Range("M34").Value = -1506.6786
$ws.Range("N34").Value = -3987.5
# row 132
$ws.Range("H132").Value = 1420.4147
$ws.Range("I132").Value = 1145.1389
$ws.Range("J132").Value = 3402.4
$ws.Range("K132").Value = 3435.4167
$ws.Range("L132").Value = 10207.2
$ws.Range("M132").Value = -905.4166999999998
$ws.Range("N132").Value = -15267.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 8
$ws.Range("H8").Value = 149
$ws.Range("I8").Value = 149
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 447
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -308
# row 14
$ws.Range("H14").Value = 802.38464
$ws.Range("I14").Value = 802.38464
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 2407.15392
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -2234.15392
# row 33
$ws.Range("H33").Value = 180
$ws.Range("I33").Value = 233.33333
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 1399.99998
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = -1116.99998
$ws.Range("N33").Value = -1166
# row 44
$ws.Range("H44").Value = 221.14285
$ws.Range("I44").Value = 269.6
$ws.Range("J44").Value = 100
$ws.Range("K44").Value = 808.8000000000001
$ws.Range("L44").Value = 300
$ws.Range("M44").Value = -410.8000000000001
$ws.Range("N44").Value = -1096
# row 64
$ws.Range("H64").Value = 2233.3333
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2233.3333
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 6699.999899999999
$ws.Range("N64").Value = -7239.999899999999
$ws.Range("M64").ClearContents()
# row 67
$ws.Range("H67").Value = 2233.3333
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2233.3333
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 6699.999899999999
$ws.Range("N67").Value = -8571.999899999999
$ws.Range("M67").ClearContents()
# row 80
$ws.Range("H80").Value = 402
$ws.Range("I80").Value = 402
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1206
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -270
# row 83
$ws.Range("H83").Value = 402
$ws.Range("I83").Value = 402
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 3618
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 1062
# row 92
$ws.Range("H92").Value = 625
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 650
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 1950
$ws.Range("M92").Value = -552
$ws.Range("N92").Value = -4446
# row 97
$ws.Range("H97").Value = 802
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 802
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2406
$ws.Range("N97").Value = -3398
$ws.Range("M97").ClearContents()
# row 98
$ws.Range("H98").Value = 762.1
$ws.Range("I98").Value = 627.8333
$ws.Range("J98").Value = 963.5
$ws.Range("K98").Value = 1883.4999
$ws.Range("L98").Value = 2890.5
$ws.Range("M98").Value = -385.4999
$ws.Range("N98").Value = -5886.5
# row 104
$ws.Range("H104").Value = 3475
$ws.Range("I104").Value = 2400
$ws.Range("J104").Value = 3628.5715
$ws.Range("K104").Value = 7200
$ws.Range("L104").Value = 10885.7145
$ws.Range("M104").Value = -4579
$ws.Range("N104").Value = -16127.7145
# row 107
$ws.Range("H107").Value = 150.57143
$ws.Range("I107").Value = 145.33333
$ws.Range("J107").Value = 154.5
$ws.Range("K107").Value = 435.99999
$ws.Range("L107").Value = 463.5
$ws.Range("M107").Value = 1484.00001
$ws.Range("N107").Value = -4303.5
# row 118
$ws.Range("H118").Value = 1635.5834
$ws.Range("I118").Value = 728.375
$ws.Range("J118").Value = 3450
$ws.Range("K118").Value = 2185.125
$ws.Range("L118").Value = 10350
$ws.Range("M118").Value = -942.125
$ws.Range("N118").Value = -12836
# row 121
$ws.Range("H121").Value = 538.8125
$ws.Range("I121").Value = 491.53845
$ws.Range("J121").Value = 743.6667
$ws.Range("K121").Value = 1474.61535
$ws.Range("L121").Value = 2231.0001
$ws.Range("M121").Value = -164.61535
$ws.Range("N121").Value = -4851.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 1160.2
$ws.Range("I46").Value = 950.5
$ws.Range("J46").Value = 1300
$ws.Range("K46").Value = 950.5
$ws.Range("L46").Value = 1300
$ws.Range("M46").Value = -762.5
$ws.Range("N46").Value = -1676
# row 61
$ws.Range("H61").Value = 1051.5
$ws.Range("I61").Value = 1051.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1051.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -849.5
$ws.Range("N61").ClearContents()
# row 113
$ws.Range("H113").Value = 1051.5
$ws.Range("I113").Value = 1051.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1051.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1118.5
$ws.Range("N113").ClearContents()
# row 130
$ws.Range("H130").Value = 59666.668
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 59666.668
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 59666.668
$ws.Range("N130").Value = -69706.66800000001
